# Re-doing global M2 module
# Applies updated M2/FX length & date-serial figures for the Top8 data comp sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - China
$ws.Range("C2").Value = 360
$ws.Range("F2").Value = 45992
$ws.Range("G2").Value = 30865
$ws.Range("H2").Value = 46055

# Row 3 - United States
$ws.Range("E3").Value = 30803
$ws.Range("F3").Value = 45992

# Row 4 - Euro Area
$ws.Range("E4").Value = 30803
$ws.Range("F4").Value = 45992
$ws.Range("G4").Value = 30865
$ws.Range("H4").Value = 46055

# Row 5 - Japan
$ws.Range("E5").Value = 30803
$ws.Range("F5").Value = 45992
$ws.Range("G5").Value = 30865
$ws.Range("H5").Value = 46055

# Row 6 - United Kingdom
$ws.Range("G6").Value = 30865
$ws.Range("H6").Value = 46055

# Row 7 - South Korea
$ws.Range("E7").Value = 30773
$ws.Range("F7").Value = 45962
$ws.Range("G7").Value = 30865
$ws.Range("H7").Value = 46055

# Row 8 - Hong Kong
$ws.Range("D8").Value = 436
$ws.Range("H8").Value = 46055

# Row 9 - Australia
$ws.Range("G9").Value = 30865
$ws.Range("H9").Value = 46055
